# feat: add 2022-Q3 data
#
# Workbook currently has two sheets: "总计" (totals) and "2022-Q2" (fund
# holdings detail for quarter 2022-Q2).
#
# This edit:
#   1. Preserves the existing "2022-Q2" detail sheet's data by duplicating
#      it into a new sheet (placed right after it).
#   2. Renames/updates the original "2022-Q2" sheet in place to become the
#      new "2022-Q3" sheet, refreshing its figures with the latest quarter
#      data (this keeps the original sheet identity/position = 2nd sheet).
#   3. Renames the duplicated sheet back to "2022-Q2" so the historical data
#      remains available as the 3rd sheet.
#   4. Updates the "总计" summary sheet: the existing summary row now
#      reports on 2022-Q3, and a new row is appended preserving the
#      2022-Q2 summary figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q2" sheet so its data is preserved
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $q2Sheet)

$oldDataSheet = $wb.Worksheets.Item(3)
$oldDataSheet.Name = "2022-Q2-old"

# ---------------------------------------------------------------------
# Step 2: turn the original sheet into the new "2022-Q3" sheet with the
# refreshed quarterly figures
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet.Name = "2022-Q3"

$q3Sheet.Range("C2").Value = "浦银安盛全球智能科技股票（QDII）A"
$q3Sheet.Range("D2").Value = "'0.25"
$q3Sheet.Range("E2").Value = "'84.65"
$q3Sheet.Range("F2").Value = "'2.53"
$q3Sheet.Range("G2").Value = "'0.0063"
$q3Sheet.Range("H2").Value = 9

$q3Sheet.Range("E3").Value = "'84.65"
$q3Sheet.Range("F3").Value = "'2.53"
$q3Sheet.Range("G3").Value = "'0.0003"
$q3Sheet.Range("H3").Value = 9

# ---------------------------------------------------------------------
# Step 3: rename the duplicated sheet back to "2022-Q2" so the old data
# keeps living under its original sheet name, now as the 3rd sheet
# ---------------------------------------------------------------------
$oldDataSheet.Name = "2022-Q2"

# ---------------------------------------------------------------------
# Step 4: update the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("B2").Value = "2022-Q3"

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.01

# match the formatting already used further up column A
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
